# Rename the visible worksheet from "walmart_test_cases" to "gaps_test_cases"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("walmart_test_cases")
$ws.Name = "gaps_test_cases"

# Update the test-case description in C9: clarify that the balance check
# should NOT succeed with a wrong gift card number.
$ws.Range("C9").Value = "User can not check gift card balance with wrong gift card number"

# Reflect the author's final selection/scroll state on the sheet.
$ws.Activate()
$ws.Range("D13").Select()
